# Update with GW7 data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 - Alisson Ramses Becker
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 3.428571428571428
$ws.Range("G2").Value = -0.1396214181607615
$ws.Range("H2").Value = 1.048225023700889
$ws.Range("I2").Value = -0.3694035501515322
$ws.Range("J2").Value = 0.3622567178255499
$ws.Range("N2").Value = 3.142857142857143

# ---------------------------------------------------------------------------
# Row 3 - Ederson Santana de Moraes
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 3.714285714285714
$ws.Range("G3").Value = 0.008449793750060434
$ws.Range("H3").Value = 0.9529078432517016
$ws.Range("I3").Value = 0.02235605289244778
$ws.Range("J3").Value = 0.4914444373324912
$ws.Range("L3").Value = "None"
$ws.Range("L3").Interior.Color = 16381413
$ws.Range("L3").Font.Color = 0
$ws.Range("N3").Value = 2.857142857142857
$ws.Range("N3").Interior.Color = 6194490
$ws.Range("N3").Font.Color = 15856113

# ---------------------------------------------------------------------------
# Row 4 - Jose Malheiro de Sa
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 3.142857142857143
$ws.Range("G4").Value = -0.125328811794971
$ws.Range("H4").Value = 1.037383927740272
$ws.Range("I4").Value = -0.3315888681207119
$ws.Range("J4").Value = 0.3757377817414266
$ws.Range("N4").Value = 3.428571428571428
$ws.Range("N4").Interior.Color = 4602842
$ws.Range("N4").Font.Color = 15856113

# ---------------------------------------------------------------------------
# Row 5 - Nick Pope
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 3.857142857142857
$ws.Range("G5").Value = 0.1955381429214488
$ws.Range("H5").Value = 0.927716703997147
$ws.Range("I5").Value = 0.5173452979975585
$ws.Range("J5").Value = 0.3117133996820917
$ws.Range("N5").Value = 3.142857142857143
$ws.Range("N5").Interior.Color = 15856114
$ws.Range("N5").Font.Color = 0

# ---------------------------------------------------------------------------
# Row 6 - Aaron Ramsdale
# ---------------------------------------------------------------------------
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 1.571428571428571
$ws.Range("G6").Value = -0.7480415953363025
$ws.Range("H6").Value = 0.8425458509342351
$ws.Range("I6").Value = -1.979132031591871
$ws.Range("J6").Value = 0.04756915711966781
$ws.Range("L6").Value = "Medium"
$ws.Range("L6").Interior.Color = 10732133
$ws.Range("L6").Font.Color = 16777215
$ws.Range("M6").Value = 4.9
$ws.Range("N6").Value = 3.142857142857143

# ---------------------------------------------------------------------------
# Row 7 - Lukasz Fabianski
# ---------------------------------------------------------------------------
$ws.Range("E7").Value = 7
$ws.Range("H7").Value = 0.9209922701055917
$ws.Range("I7").Value = -3.485942738065964
$ws.Range("J7").Value = 0.006523834583348492
$ws.Range("K7").Value = $true
$ws.Range("K7").Interior.Color = 2263842
$ws.Range("K7").Font.Color = 16777215
$ws.Range("L7").Value = "Very large"
$ws.Range("L7").Interior.Color = 4491810
$ws.Range("L7").Font.Color = 16777215
$ws.Range("N7").Value = 3.142857142857143
$ws.Range("N7").Interior.Color = 15856114
$ws.Range("N7").Font.Color = 0

# ---------------------------------------------------------------------------
# Row 8 - Emiliano Martinez Romero
# ---------------------------------------------------------------------------
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 2.571428571428572
$ws.Range("G8").Value = -0.2251445479836814
$ws.Range("H8").Value = 0.9845380540342235
$ws.Range("I8").Value = -0.5956764830068697
$ws.Range("J8").Value = 0.2865795267040163
$ws.Range("L8").Value = "Small"
$ws.Range("L8").Interior.Color = 13228184
$ws.Range("L8").Font.Color = 16777215
$ws.Range("N8").Value = 3.142857142857143

# ---------------------------------------------------------------------------
# Row 9 - Hugo Lloris
# ---------------------------------------------------------------------------
$ws.Range("E9").Value = 7
$ws.Range("H9").Value = 0.9293350223348444
$ws.Range("I9").Value = -3.398015420648054
$ws.Range("J9").Value = 0.007265699366264177
$ws.Range("K9").Value = $true
$ws.Range("K9").Interior.Color = 2263842
$ws.Range("K9").Font.Color = 16777215
$ws.Range("L9").Value = "Very large"
$ws.Range("L9").Interior.Color = 4491810
$ws.Range("L9").Font.Color = 16777215
$ws.Range("M9").Value = 4.8
$ws.Range("N9").Value = 3.142857142857143
$ws.Range("N9").Interior.Color = 15856114
$ws.Range("N9").Font.Color = 0

# ---------------------------------------------------------------------------
# Row 10 - David Raya Martin
# ---------------------------------------------------------------------------
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = -0.4578014486081555
$ws.Range("H10").Value = 0.912513947433491
$ws.Range("I10").Value = -1.211228782862296
$ws.Range("J10").Value = 0.1356721871696142
$ws.Range("L10").Value = "Small"
$ws.Range("L10").Interior.Color = 13228184
$ws.Range("L10").Font.Color = 16777215
$ws.Range("N10").Value = 3.142857142857143

# ---------------------------------------------------------------------------
# Row 11 - Jordan Pickford
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 1.571428571428571
$ws.Range("G11").Value = -0.5671818103917637
$ws.Range("H11").Value = 0.8467666581600165
$ws.Range("I11").Value = -1.500622018455997
$ws.Range("J11").Value = 0.09206230105601441
$ws.Range("L11").Value = "Medium"
$ws.Range("L11").Interior.Color = 10732133
$ws.Range("L11").Font.Color = 16777215
$ws.Range("N11").Value = 2.857142857142857
$ws.Range("N11").Interior.Color = 6194490
$ws.Range("N11").Font.Color = 15856113
